$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "26/08/2016"
$ws.Range("B7:J7").Value = "done"

$ws.Range("B7").Select()
